# Weekly update: insert a new price record for
# "Macroferia Regional de Talca - Zanahoria" as the new row 335,
# pushing the existing rows 335-359 down to 336-360.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 335 (shifts 335:359 -> 336:360).
$ws.Rows.Item(335).Insert()

# Populate the new row 335 with this week's data.
$ws.Cells.Item(335, 1).Value  = 5
$ws.Cells.Item(335, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(335, 3).Value  = "Maule"
$ws.Cells.Item(335, 4).Value  = 44746
$ws.Cells.Item(335, 5).Value  = 7
$ws.Cells.Item(335, 6).Value  = 100114013
$ws.Cells.Item(335, 7).Value  = "Zanahoria"
$ws.Cells.Item(335, 8).Value  = "Sin especificar"
$ws.Cells.Item(335, 9).Value  = "Primera"
$ws.Cells.Item(335, 10).Value = 500
$ws.Cells.Item(335, 11).Value = 6000
$ws.Cells.Item(335, 12).Value = 6000
$ws.Cells.Item(335, 13).Value = 6000
$ws.Cells.Item(335, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(335, 15).Value = "Región de Ñuble"
$ws.Cells.Item(335, 16).Value = 300
$ws.Cells.Item(335, 17).Value = 20
$ws.Cells.Item(335, 18).Value = "Hortaliza"

# Give the date cell the same date style used by the rest of column D.
$ws.Cells.Item(335, 4).Style = $ws.Cells.Item(336, 4).Style
